$wb = $excel.ActiveWorkbook

# Remove the "Metadata" sheet entirely.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Metadata").Delete()
$excel.DisplayAlerts = $true

# The remaining "Reference" sheet becomes the single sheet, renamed to "Sheet1".
$ws = $wb.Worksheets.Item("Reference")
$ws.Name = "Sheet1"

# Move the on-sheet selection to A32.
$ws.Activate()
$ws.Range("A32").Select()
